$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.279.47"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.620.32"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.06"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.07"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "3.609.33"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -4.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.83"
$ws.Range("E11").Value = "  +23.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.603"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.21"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "4.203.34"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "666.90"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.87"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "3.616.00"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("D19").Value = "70.282.79"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.73"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.932"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.06"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.53"
$ws.Range("E25").Value = "  -4.44%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.60"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.95"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("E33").Value = "  -6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.36"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "576.36"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.02"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.10"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "3.573.15"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.343"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.43"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("D46").Value = "0.0₃0734"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -4.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.83"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.96"
$ws.Range("E51").Value = "  +2.67%  "
